$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

# Columns whose literal text would otherwise be auto-converted by Excel
# (dates, and numeric-looking strings like "00") need to be forced to
# text and then have their formatting cleared again so no explicit
# cell style ends up attached (matching the rest of the data rows).
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value2 = "2024-01-04"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value2 = "11:55:27"

$ws.Range("C$row").Value2 = "Thursday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value2 = "00"
$ws.Range("D$row").ClearFormats()

$ws.Range("E$row").Value2 = 140262
$ws.Range("F$row").Value2 = 142897
$ws.Range("G$row").Value2 = 171950
$ws.Range("H$row").Value2 = 146990
$ws.Range("I$row").Value2 = -1
$ws.Range("J$row").Value2 = 117482
$ws.Range("K$row").Value2 = 224053
$ws.Range("L$row").Value2 = 248087
$ws.Range("M$row").Value2 = 184304
$ws.Range("N$row").Value2 = 109999
$ws.Range("O$row").Value2 = 40287
$ws.Range("P$row").Value2 = 30837
$ws.Range("Q$row").Value2 = 72298
$ws.Range("R$row").Value2 = -1
$ws.Range("S$row").Value2 = 41240
$ws.Range("T$row").Value2 = -1
